# Project_Notebook.xlsx - revision bump edit
# Commit: "UpLoaded on 9:10am 17May23 by Willy Added a 1sec delay between
#          rewriting clipboard value to remove read error leading to the
#          halt of Everest."
#
# The only real data change the author made inside the workbook itself is
# on the "File" sheet: the input full file path (FullFileName_Input,
# File!B6) had its revision number reset from 388 down to 0. Every other
# cell on that sheet (FileName, File_Name_No_Ext, CurrentRevision,
# NewRevision, the New_DML_*/New_Folder_* paths, etc.) is a formula that
# is a function of B6, so updating B6 ripples through automatically.
#
# (The Date sheet's TODAY()-driven cells move forward one day because the
# workbook was recalculated the next morning - that happens automatically
# on recalculation against the correct wall-clock date and needs no
# explicit write here.)

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("File")

$ws.Range("B6").Value = "C:\Temp\Everest\Document masterlist\CJA\CJA Document Masterlist Everest Rev 0.xlsx"
